$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text could be misread as a number by Excel's type inference;
# these get a temporary Text number format so the literal string is preserved,
# then the style is reset back to Normal so no stray formatting is introduced.
$textSafeCells = @(
    @{Cell='D5'; Value='219.51'},
    @{Cell='D6'; Value='0.5168'},
    @{Cell='D8'; Value='0.06474'},
    @{Cell='D9'; Value='0.2580'},
    @{Cell='D10'; Value='20.00'},
    @{Cell='D11'; Value='0.07684'},
    @{Cell='D12'; Value='4.353'},
    @{Cell='D15'; Value='0.5584'},
    @{Cell='D17'; Value='64.92'},
    @{Cell='D19'; Value='1.009'},
    @{Cell='D20'; Value='211.02'},
    @{Cell='D21'; Value='4.447'},
    @{Cell='D22'; Value='10.11'},
    @{Cell='D23'; Value='5.899'},
    @{Cell='D25'; Value='143.21'},
    @{Cell='D26'; Value='1.719'},
    @{Cell='D28'; Value='7.004'},
    @{Cell='D29'; Value='15.80'},
    @{Cell='D30'; Value='0.05227'},
    @{Cell='D32'; Value='3.354'},
    @{Cell='D33'; Value='3.217'},
    @{Cell='D34'; Value='1.586'},
    @{Cell='D35'; Value='2.762'},
    @{Cell='D36'; Value='2.379'},
    @{Cell='D37'; Value='0.9251'},
    @{Cell='D38'; Value='0.5738'},
    @{Cell='D40'; Value='0.01594'},
    @{Cell='D42'; Value='0.8392'},
    @{Cell='D43'; Value='5.647'},
    @{Cell='D44'; Value='100.03'},
    @{Cell='D47'; Value='0.4497'},
    @{Cell='D48'; Value='55.87'},
    @{Cell='D50'; Value='7.950'}
)

foreach ($item in $textSafeCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
    $rng.Style = "Normal"
}

# Plain text updates (not at risk of numeric auto-conversion)
$plainCells = @(
    @{Cell='D2'; Value='26.509.34'},
    @{Cell='E2'; Value='  -2.84%  '},
    @{Cell='D3'; Value='1.671.19'},
    @{Cell='E3'; Value='  -2.07%  '},
    @{Cell='E4'; Value='  +0.45%  '},
    @{Cell='E5'; Value='  -1.72%  '},
    @{Cell='E6'; Value='  -2.59%  '},
    @{Cell='E7'; Value='  +0.40%  '},
    @{Cell='E8'; Value='  -1.68%  '},
    @{Cell='E9'; Value='  -2.71%  '},
    @{Cell='E10'; Value='  -4.03%  '},
    @{Cell='E11'; Value='  +0.60%  '},
    @{Cell='B12'; Value='Polkadot'},
    @{Cell='C12'; Value='https://coinranking.com/coin/25W7FG7om+polkadot-dot'},
    @{Cell='E12'; Value='  -4.88%  '},
    @{Cell='B13'; Value='WrappedliquidstakedEther2.0'},
    @{Cell='C13'; Value='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'},
    @{Cell='D13'; Value='1.901.25'},
    @{Cell='E13'; Value='  -2.12%  '},
    @{Cell='B14'; Value='WrappedEther'},
    @{Cell='C14'; Value='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'},
    @{Cell='D14'; Value='1.664.50'},
    @{Cell='E14'; Value='  -1.94%  '},
    @{Cell='E15'; Value='  -2.51%  '},
    @{Cell='D16'; Value='0.0₅8036'},
    @{Cell='E16'; Value='  -1.74%  '},
    @{Cell='E17'; Value='  -3.87%  '},
    @{Cell='D18'; Value='26.544.79'},
    @{Cell='E18'; Value='  -2.71%  '},
    @{Cell='E19'; Value='  +0.59%  '},
    @{Cell='E20'; Value='  -2.33%  '},
    @{Cell='E21'; Value='  -4.70%  '},
    @{Cell='E22'; Value='  -2.89%  '},
    @{Cell='E23'; Value='  -1.14%  '},
    @{Cell='E24'; Value='  +0.38%  '},
    @{Cell='E25'; Value='  +0.98%  '},
    @{Cell='E26'; Value='  -2.11%  '},
    @{Cell='E27'; Value='  -3.79%  '},
    @{Cell='E28'; Value='  -3.48%  '},
    @{Cell='E29'; Value='  -3.18%  '},
    @{Cell='E30'; Value='  -2.94%  '},
    @{Cell='E31'; Value='  -1.99%  '},
    @{Cell='E32'; Value='  -4.19%  '},
    @{Cell='E33'; Value='  -5.89%  '},
    @{Cell='E34'; Value='  -3.19%  '},
    @{Cell='E35'; Value='  -4.06%  '},
    @{Cell='E36'; Value='  -1.87%  '},
    @{Cell='E37'; Value='  -2.24%  '},
    @{Cell='E38'; Value='  -1.84%  '},
    @{Cell='D39'; Value='1.159.02'},
    @{Cell='E39'; Value='  +11.15%  '},
    @{Cell='E40'; Value='  -2.44%  '},
    @{Cell='E41'; Value='  +0.39%  '},
    @{Cell='E42'; Value='  +0.05%  '},
    @{Cell='E43'; Value='  -3.66%  '},
    @{Cell='E44'; Value='  -0.93%  '},
    @{Cell='D45'; Value='1.810.17'},
    @{Cell='E45'; Value='  -2.15%  '},
    @{Cell='E46'; Value='  -2.32%  '},
    @{Cell='E47'; Value='  -0.12%  '},
    @{Cell='E48'; Value='  -3.71%  '},
    @{Cell='E49'; Value='  -0.18%  '},
    @{Cell='E50'; Value='  -1.83%  '},
    @{Cell='E51'; Value='  -1.94%  '}
)

foreach ($item in $plainCells) {
    $ws.Range($item.Cell).Value = $item.Value
}
